$d = $word.ActiveDocument

# Update the date/title paragraph
$d.Paragraphs.Item(1).Range.Text = "2025-08-30 Saturday"

# Update the table cells (positional addressing avoids ambiguity
# since some new values equal other cells old values)
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "30÷5="
$t.Cell(1,2).Range.Text = "96÷9="
$t.Cell(1,3).Range.Text = "75÷7="
$t.Cell(1,4).Range.Text = "28÷6="
$t.Cell(1,5).Range.Text = "82÷2="
$t.Cell(5,1).Range.Text = "63÷9="
$t.Cell(5,2).Range.Text = "55÷5="
$t.Cell(5,3).Range.Text = "39÷8="
$t.Cell(5,4).Range.Text = "64÷4="
$t.Cell(5,5).Range.Text = "11÷5="
$t.Cell(9,1).Range.Text = "71÷7="
$t.Cell(9,2).Range.Text = "47÷3="
$t.Cell(9,3).Range.Text = "62÷4="
$t.Cell(9,4).Range.Text = "20÷2="
$t.Cell(9,5).Range.Text = "29÷5="
$t.Cell(13,1).Range.Text = "62÷4="
$t.Cell(13,2).Range.Text = "31÷7="
$t.Cell(13,3).Range.Text = "77÷7="
$t.Cell(13,4).Range.Text = "31÷5="
$t.Cell(13,5).Range.Text = "11÷6="
$t.Cell(17,1).Range.Text = "53÷7="
$t.Cell(17,2).Range.Text = "84÷3="
$t.Cell(17,3).Range.Text = "49÷8="
$t.Cell(17,4).Range.Text = "90÷6="
$t.Cell(17,5).Range.Text = "82÷5="
